# Generate Report for Handoff
# The "b.md" file has moved from "Handed back: in sync with en-US" to
# "Ready for handoff" for both locales. A new (out-of-date) handback file
# was received for zh-cn, and for de-de the handback turned out to be stale,
# so Content Duplicate flips to False and an Error Detail message is recorded.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/75be731a63683d49190e51e40793af59acf1723b/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a724c00367cc6a8703b4ef610727dc28f83eb2da/e2e/b.md."

# ---------------------------------------------------------------------------
# Overview sheet: row 3 is the b.md file. Its zh-cn / de-de status columns and
# the "Latest HO Xliff Generate Date" column move forward.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-19 04:33:07"

# ---------------------------------------------------------------------------
# zh-cn sheet: row 3 (b.md) gets a new, out-of-date handback file.
# (leading "'" forces the "True"/"False" text to be stored as text, not a
#  real boolean, matching the rest of the sheet's True/False text columns)
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("F3").Style = "Normal"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-19 04:32:58"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = 39.14

# ---------------------------------------------------------------------------
# de-de sheet: row 3 (b.md) gets a new, out-of-date handback file.
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "'False"
$dede.Range("F3").Style = "Normal"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-19 04:33:07"
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = 39.14
